# "fix source documents that used h1 for normal text": the "None"
# paragraph under "Prerequisite tasks" is wrongly styled as Heading 1 -
# it should just be normal body text. Fixing the style also removes the
# stray sz=22/szCs=22 direct-formatting override that came along with
# the heading, and - because the bookmark that used to sit on this
# heading paragraph is no longer needed once it's a plain paragraph -
# that bookmark is removed too (which causes Word to renumber the
# bookmark ids that follow it).

$d = $word.ActiveDocument

# Locate the paragraph index that holds the "_aqdfz55armzg" bookmark
# (the "None" line) instead of hard-coding an index.
$bm = $d.Bookmarks.Item("_aqdfz55armzg")
$bmStart = $bm.Start
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Start -le $bmStart -and $candidate.Range.End -gt $bmStart) {
        $targetIndex = $i
    }
}

$p = $d.Paragraphs.Item($targetIndex)

# 1. Re-style it as Normal (plain body text) instead of Heading 1.
$p.Style = "Normal"

# 2. Clear the leftover direct character formatting (sz/szCs=22) and the
#    stray empty run trailing it by deleting the run text and retyping
#    it fresh - the new run picks up the (now Normal) paragraph's plain
#    formatting instead of carrying over the old heading run props.
$p = $d.Paragraphs.Item($targetIndex)
$textRng = $d.Range($p.Range.Start, $p.Range.End - 1)
$textRng.Delete()
$insertRng = $d.Range($p.Range.Start, $p.Range.Start)
$insertRng.InsertAfter("None")

# 3. The bookmark that lived on the old Heading-1 paragraph is no longer
#    needed - drop it. Word renumbers the remaining bookmark ids
#    (_8e3vikzea5p, _g4mxm0ixbdvd, _3r7k4knw2b7d) down by one
#    automatically.
$d.Bookmarks.Item("_aqdfz55armzg").Delete()
